{"js": "// Update the last \"CHOWCHOW\" purchase record (THU Nov 09 block):\n//   - timestamp        \" 10:56:55 PST 2017\" -> \" 10:36:44 PST 2017\"\n//   - Rate              \"- 53\"      -> \"- 12\"\n//   - Total Price       \"- 2756.0\"  -> \"- 624.0\"\n//   - Amount balance    \"- 2756.0\"  -> \"- 624.0\"\n\nconst body = context.document.body;\n\n// 1) Timestamp: unique string in the document, safe to match directly.\nconst tsResults = body.search(\" 10:56:55 PST 2017\", { matchCase: true });\ntsResults.load(\"items\");\nawait context.sync();\nfor (const r of tsResults.items) {\n  r.insertText(\" 10:36:44 PST 2017\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Rate line: \"- 53\" alone is ambiguous as a plain substring (it also sits\n//    inside \"- 5360\" / \"- 5360.0\" elsewhere in the document), so require a\n//    whole-word match to land on just the standalone \"- 53\" value.\nconst rateResults = body.search(\"- 53\", { matchCase: true, matchWholeWord: true });\nrateResults.load(\"items\");\nawait context.sync();\nfor (const r of rateResults.items) {\n  r.insertText(\"- 12\", \"Replace\");\n}\nawait context.sync();\n\n// 3) \"- 2756.0\" occurs exactly twice in the document (Total Price line and\n//    the bold Amount balance line of the same record) - both change to\n//    \"- 624.0\".\nconst totalResults = body.search(\"- 2756.0\", { matchCase: true, matchWholeWord: true });\ntotalResults.load(\"items\");\nawait context.sync();\nfor (const r of totalResults.items) {\n  r.insertText(\"- 624.0\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the last \"CHOWCHOW\" purchase record (THU Nov 09 block):\n#   - timestamp        \" 10:56:55 PST 2017\" -> \" 10:36:44 PST 2017\"\n#   - Rate              \"- 53\"      -> \"- 12\"\n#   - Total Price       \"- 2756.0\"  -> \"- 624.0\"\n#   - Amount balance    \"- 2756.0\"  -> \"- 624.0\"\n\n$d = $word.ActiveDocument\n\n# 1) Timestamp: unique string in the document, safe to match directly.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\" 10:56:55 PST 2017\", $true, $false, $false, $false, $false, $true, 1, $false, \" 10:36:44 PST 2017\", 1)\n\n# 2) Rate line: \"- 53\" alone is ambiguous as a plain substring (it also sits\n#    inside \"- 5360\" / \"- 5360.0\" elsewhere in the document), so require a\n#    whole-word match to land on just the standalone \"- 53\" value.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"- 53\", $true, $true, $false, $false, $false, $true, 1, $false, \"- 12\", 1)\n\n# 3) \"- 2756.0\" occurs exactly twice in the document (Total Price line and\n#    the bold Amount balance line of the same record) - both change to\n#    \"- 624.0\". ReplaceAll (2) picks up both matches in one pass.\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Execute(\"- 2756.0\", $true, $true, $false, $false, $false, $true, 1, $false, \"- 624.0\", 2)\n"}
